$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the column-A style (bold, centered, bordered) down to the new rows (23-25)
# that did not exist before, so A23:A25 match the formatting of A2:A22.
$ws.Range("A22").Copy()
$ws.Range("A23:A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "model_11_0_0"
$ws.Range("B2").Value = -0.06816068755473803
$ws.Range("C2").Value = -0.1145408450526688
$ws.Range("D2").Value = 0.01029391316125328
$ws.Range("E2").Value = 0.007358617732927009
$ws.Range("F2").Value = 1.182139158248901
$ws.Range("G2").Value = 1.336656451225281
$ws.Range("H2").Value = 1.93499755859375
$ws.Range("I2").Value = 1.618228673934937

$ws.Range("A3").Value = "model_11_0_1"
$ws.Range("B3").Value = -0.04800877986950769
$ws.Range("C3").Value = -0.1081633196274912
$ws.Range("D3").Value = 0.05274911024270468
$ws.Range("E3").Value = 0.03380314409299034
$ws.Range("F3").Value = 1.159836888313293
$ws.Range("G3").Value = 1.329007983207703
$ws.Range("H3").Value = 1.85199236869812
$ws.Range("I3").Value = 1.575118184089661

$ws.Range("A4").Value = "model_11_0_2"
$ws.Range("B4").Value = 0.04625563942859168
$ws.Range("C4").Value = -0.0008710964236273711
$ws.Range("D4").Value = 0.0389028746495631
$ws.Range("E4").Value = 0.06777532426969946
$ws.Range("F4").Value = 1.055513978004456
$ws.Range("G4").Value = 1.200333714485168
$ws.Range("H4").Value = 1.879063606262207
$ws.Range("I4").Value = 1.519735932350159

$ws.Range("A5").Value = "model_11_0_3"
$ws.Range("B5").Value = 0.08832561532274152
$ws.Range("C5").Value = 0.03279447349524423
$ws.Range("D5").Value = 0.05104451970321155
$ws.Range("E5").Value = 0.08773935075752892
$ws.Range("F5").Value = 1.008954882621765
$ws.Range("G5").Value = 1.159958839416504
$ws.Range("H5").Value = 1.855325102806091
$ws.Range("I5").Value = 1.487190008163452

$ws.Range("A6").Value = "model_11_0_4"
$ws.Range("B6").Value = 0.3156824169224492
$ws.Range("C6").Value = 0.1080818448311355
$ws.Range("D6").Value = 0.2932552635074922
$ws.Range("E6").Value = 0.2537589461622238
$ws.Range("F6").Value = 0.7573379278182983
$ws.Range("G6").Value = 1.06966757774353
$ws.Range("H6").Value = 1.381773114204407
$ws.Range("I6").Value = 1.216540813446045

$ws.Range("A7").Value = "model_11_0_5"
$ws.Range("B7").Value = 0.331640229170591
$ws.Range("C7").Value = 0.1146681486996772
$ws.Range("D7").Value = 0.2942122728580427
$ws.Range("E7").Value = 0.2568642014313111
$ws.Range("F7").Value = 0.7396773099899292
$ws.Range("G7").Value = 1.061768651008606
$ws.Range("H7").Value = 1.379902124404907
$ws.Range("I7").Value = 1.211478590965271

$ws.Range("A8").Value = "model_11_0_8"
$ws.Range("B8").Value = 0.3450483649766604
$ws.Range("C8").Value = 0.1211862141628894
$ws.Range("D8").Value = 0.2509982088256177
$ws.Range("E8").Value = 0.2350138551476394
$ws.Range("F8").Value = 0.7248384356498718
$ws.Range("G8").Value = 1.053951621055603
$ws.Range("H8").Value = 1.464390873908997
$ws.Range("I8").Value = 1.24709951877594

$ws.Range("A9").Value = "model_11_0_6"
$ws.Range("B9").Value = 0.3453398061520659
$ws.Range("C9").Value = 0.1287726956855113
$ws.Range("D9").Value = 0.2605987210437427
$ws.Range("E9").Value = 0.2433868260784946
$ws.Range("F9").Value = 0.7245159149169922
$ws.Range("G9").Value = 1.044853329658508
$ws.Range("H9").Value = 1.445620775222778
$ws.Range("I9").Value = 1.233449697494507

$ws.Range("A10").Value = "model_11_0_7"
$ws.Range("B10").Value = 0.3470543219751281
$ws.Range("C10").Value = 0.1256730746258826
$ws.Range("D10").Value = 0.2646725653375472
$ws.Range("E10").Value = 0.2444787706688515
$ws.Range("F10").Value = 0.7226185202598572
$ws.Range("G10").Value = 1.04857063293457
$ws.Range("H10").Value = 1.437655925750732
$ws.Range("I10").Value = 1.231669545173645

$ws.Range("A11").Value = "model_11_0_10"
$ws.Range("B11").Value = 0.3490945284164856
$ws.Range("C11").Value = 0.09629133287235891
$ws.Range("D11").Value = 0.2659729499025371
$ws.Range("E11").Value = 0.2337695074778183
$ws.Range("F11").Value = 0.7203605771064758
$ws.Range("G11").Value = 1.083807706832886
$ws.Range("H11").Value = 1.435113430023193
$ws.Range("I11").Value = 1.249127984046936

$ws.Range("A12").Value = "model_11_0_11"
$ws.Range("B12").Value = 0.3491885339425902
$ws.Range("C12").Value = 0.09655965688814394
$ws.Range("D12").Value = 0.2661077711725272
$ws.Range("E12").Value = 0.2339500895612397
$ws.Range("F12").Value = 0.720256507396698
$ws.Range("G12").Value = 1.083486080169678
$ws.Range("H12").Value = 1.434849858283997
$ws.Range("I12").Value = 1.248833656311035

$ws.Range("A13").Value = "model_11_0_13"
$ws.Range("B13").Value = 0.349212249606887
$ws.Range("C13").Value = 0.09549771623655212
$ws.Range("D13").Value = 0.2668183630552963
$ws.Range("E13").Value = 0.233937463362997
$ws.Range("F13").Value = 0.7202302217483521
$ws.Range("G13").Value = 1.084759593009949
$ws.Range("H13").Value = 1.433460593223572
$ws.Range("I13").Value = 1.248854279518127

$ws.Range("A14").Value = "model_11_0_12"
$ws.Range("B14").Value = 0.3492239351016028
$ws.Range("C14").Value = 0.09656495287040912
$ws.Range("D14").Value = 0.266207479273331
$ws.Range("E14").Value = 0.2340084496531647
$ws.Range("F14").Value = 0.7202172875404358
$ws.Range("G14").Value = 1.083479642868042
$ws.Range("H14").Value = 1.434654951095581
$ws.Range("I14").Value = 1.248738527297974

$ws.Range("A15").Value = "model_11_0_14"
$ws.Range("B15").Value = 0.349372772568099
$ws.Range("C15").Value = 0.09440801020418665
$ws.Range("D15").Value = 0.268158339668005
$ws.Range("E15").Value = 0.2342693739225719
$ws.Range("F15").Value = 0.7200526595115662
$ws.Range("G15").Value = 1.086066484451294
$ws.Range("H15").Value = 1.430840730667114
$ws.Range("I15").Value = 1.248313188552856

$ws.Range("A16").Value = "model_11_0_15"
$ws.Range("B16").Value = 0.3494467356073563
$ws.Range("C16").Value = 0.09406023644540007
$ws.Range("D16").Value = 0.2686885312312213
$ws.Range("E16").Value = 0.2344330967225927
$ws.Range("F16").Value = 0.7199707627296448
$ws.Range("G16").Value = 1.086483478546143
$ws.Range("H16").Value = 1.42980420589447
$ws.Range("I16").Value = 1.248046278953552

$ws.Range("A17").Value = "model_11_0_17"
$ws.Range("B17").Value = 0.3494523917728968
$ws.Range("C17").Value = 0.09410319052176697
$ws.Range("D17").Value = 0.2686783949848724
$ws.Range("E17").Value = 0.2344441734623272
$ws.Range("F17").Value = 0.7199645042419434
$ws.Range("G17").Value = 1.086432099342346
$ws.Range("H17").Value = 1.429824113845825
$ws.Range("I17").Value = 1.24802827835083

$ws.Range("A18").Value = "model_11_0_16"
$ws.Range("B18").Value = 0.3494529252963493
$ws.Range("C18").Value = 0.09407053642250829
$ws.Range("D18").Value = 0.2687033971542466
$ws.Range("E18").Value = 0.2344455689428081
$ws.Range("F18").Value = 0.7199639678001404
$ws.Range("G18").Value = 1.086471080780029
$ws.Range("H18").Value = 1.42977511882782
$ws.Range("I18").Value = 1.248026013374329

$ws.Range("A19").Value = "model_11_0_18"
$ws.Range("B19").Value = 0.3494662924573493
$ws.Range("C19").Value = 0.09405359459439744
$ws.Range("D19").Value = 0.2687654045540175
$ws.Range("E19").Value = 0.2344739850467156
$ws.Range("F19").Value = 0.7199491858482361
$ws.Range("G19").Value = 1.086491465568542
$ws.Range("H19").Value = 1.429654002189636
$ws.Range("I19").Value = 1.247979640960693

$ws.Range("A20").Value = "model_11_0_19"
$ws.Range("B20").Value = 0.3494668663482472
$ws.Range("C20").Value = 0.09405382462940537
$ws.Range("D20").Value = 0.2687674557160448
$ws.Range("E20").Value = 0.2344751998668135
$ws.Range("F20").Value = 0.7199484705924988
$ws.Range("G20").Value = 1.086491227149963
$ws.Range("H20").Value = 1.429649949073792
$ws.Range("I20").Value = 1.247977614402771

$ws.Range("A21").Value = "model_11_0_22"
$ws.Range("B21").Value = 0.3494677884409869
$ws.Range("C21").Value = 0.09405241784367024
$ws.Range("D21").Value = 0.2687719877075746
$ws.Range("E21").Value = 0.2344772362615105
$ws.Range("F21").Value = 0.7199474573135376
$ws.Range("G21").Value = 1.086492896080017
$ws.Range("H21").Value = 1.429641127586365
$ws.Range("I21").Value = 1.247974395751953

$ws.Range("A22").Value = "model_11_0_20"
$ws.Range("B22").Value = 0.3494677884409869
$ws.Range("C22").Value = 0.09405241784367024
$ws.Range("D22").Value = 0.2687719877075746
$ws.Range("E22").Value = 0.2344772362615105
$ws.Range("F22").Value = 0.7199474573135376
$ws.Range("G22").Value = 1.086492896080017
$ws.Range("H22").Value = 1.429641127586365
$ws.Range("I22").Value = 1.247974395751953

$ws.Range("A23").Value = "model_11_0_21"
$ws.Range("B23").Value = 0.3494677884409869
$ws.Range("C23").Value = 0.09405241784367024
$ws.Range("D23").Value = 0.2687719877075746
$ws.Range("E23").Value = 0.2344772362615105
$ws.Range("F23").Value = 0.7199474573135376
$ws.Range("G23").Value = 1.086492896080017
$ws.Range("H23").Value = 1.429641127586365
$ws.Range("I23").Value = 1.247974395751953

$ws.Range("A24").Value = "model_11_0_23"
$ws.Range("B24").Value = 0.3494677884409869
$ws.Range("C24").Value = 0.09405241784367024
$ws.Range("D24").Value = 0.2687719877075746
$ws.Range("E24").Value = 0.2344772362615105
$ws.Range("F24").Value = 0.7199474573135376
$ws.Range("G24").Value = 1.086492896080017
$ws.Range("H24").Value = 1.429641127586365
$ws.Range("I24").Value = 1.247974395751953

$ws.Range("A25").Value = "model_11_0_9"
$ws.Range("B25").Value = 0.3497968985498837
$ws.Range("C25").Value = 0.09807191410583804
$ws.Range("D25").Value = 0.2678309602037704
$ws.Range("E25").Value = 0.2355115197771723
$ws.Range("F25").Value = 0.7195832133293152
$ws.Range("G25").Value = 1.081672430038452
$ws.Range("H25").Value = 1.431480884552002
$ws.Range("I25").Value = 1.246288180351257
